# Apply the "Updates in GPT versions" edit:
#  - Remove the last 7 slides (old slides 9-15 / sldId 264-270)
#  - Re-theme the remaining 8 slides from a "Numpy" training deck into an
#    "IKEA drone delivery" business-plan deck (title + 5 bullet points each).
#
# Note: when a TextRange.Text assignment produces the same paragraph count
# as the text currently in the shape, the host tries to preserve per-run
# formatting by diffing old vs new paragraph text, which can split a
# paragraph's text into two runs at the common-prefix boundary (e.g. "- "
# vs the rest). Setting a throwaway single-line value first guarantees the
# paragraph count never matches on the real assignment, so each paragraph
# always ends up as a single clean run.

$p = $ppt.ActivePresentation
$rsquo = [char]0x2019

function Set-BodyText($shape, [string]$text) {
    $shape.TextFrame.TextRange.Text = "-"
    $shape.TextFrame.TextRange.Text = $text
}

# 1) Delete the trailing slides (iterate backwards so indices stay valid)
for ($i = $p.Slides.Count; $i -ge 9; $i--) {
    $p.Slides.Item($i).Delete()
}

# 2) Rewrite title + body text for the 8 remaining slides

$s = $p.Slides.Item(1)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Market Overview"
$txt = (
    "- Rapid growth in e-commerce and delivery services  `r" +
    "- Increasing consumer demand for faster delivery options  `r" +
    "- Adoption of drone technology in logistics and retail  `r" +
    "- Competitive landscape including other retailers and delivery companies  `r" +
    "- Regulatory environment impacting drone usage"
)
Set-BodyText $s.Shapes.Item(2) $txt

$s = $p.Slides.Item(2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Target Audience"
$txt = (
    "- Primary audience: Tech-savvy millennials and Gen Z  `r" +
    "- Secondary audience: Busy professionals and families  `r" +
    "- Geographic focus: Urban and suburban areas with high order volumes  `r" +
    "- Behavioral insights: Preference for convenience and speed in shopping  `r" +
    "- Environmental concerns driving interest in sustainable delivery options"
)
Set-BodyText $s.Shapes.Item(2) $txt

$s = $p.Slides.Item(3)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Value Proposition"
$txt = (
    "- Speed: Deliver products within 30 minutes of order placement  `r" +
    "- Convenience: Doorstep delivery without human contact  `r" +
    "- Sustainability: Lower carbon footprint compared to traditional delivery  `r" +
    "- Accessibility: Reach customers in hard-to-access areas  `r" +
    "- Innovation: Strength and credibility of IKEA brand in offering cutting-edge solutions"
)
Set-BodyText $s.Shapes.Item(2) $txt

$s = $p.Slides.Item(4)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Operational Strategy"
$txt = (
    "- Partner with local logistics providers for infrastructure  `r" +
    "- Develop a network of drone landing zones for efficient pickups and drop-offs  `r" +
    "- Implement a robust technology stack for order processing and tracking  `r" +
    "- Ensure compliance with aviation regulations and safety standards  `r" +
    "- Train staff and customers on how to use the drone delivery service"
)
Set-BodyText $s.Shapes.Item(2) $txt

$s = $p.Slides.Item(5)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Marketing Approach"
$txt = (
    "- Multi-channel marketing campaign (online and offline)  `r" +
    "- Collaborate with influencers and tech bloggers to drive awareness  `r" +
    "- Demonstrate drone delivery through live events and demonstrations  `r" +
    "- Use targeted digital advertising to reach key demographics  `r" +
    "- Leverage social media to highlight customer testimonials and success stories"
)
Set-BodyText $s.Shapes.Item(2) $txt

$s = $p.Slides.Item(6)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Key Partnerships"
$txt = (
    "- Collaborate with drone manufacturers for tech innovation  `r" +
    "- Partner with local authorities to navigate regulations  `r" +
    "- Tie-up with delivery service platforms for last-mile logistics  `r" +
    "- Engage with sustainability organizations to enhance brand image  `r" +
    "- Work with data analytics firms to optimize operations"
)
Set-BodyText $s.Shapes.Item(2) $txt

$s = $p.Slides.Item(7)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Financial Projections"
$txt = (
    "- Initial investment required for drone fleet and technology  `r" +
    "- Expected return on investment (ROI) within the first 3 years  `r" +
    "- Revenue growth through increased sales and customer acquisition  `r" +
    "- Cost savings from decreased labor and fuel expenses over time  `r" +
    "- Long-term financial benefits from strengthening IKEA" + $rsquo + "s market position"
)
Set-BodyText $s.Shapes.Item(2) $txt

$s = $p.Slides.Item(8)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Implementation Timeline"
$txt = (
    "- Phase 1: Research & Development (0-6 months)  `r" +
    "- Phase 2: Pilot program launch in select cities (6-12 months)  `r" +
    "- Phase 3: Full-scale rollout across urban areas (12-24 months)  `r" +
    "- Phase 4: Evaluation and iterative enhancements (24+ months)  `r" +
    "- Continuous monitoring of market feedback and operational challenges"
)
Set-BodyText $s.Shapes.Item(2) $txt
